# WG_Number_excel_table.xlsx — renumber N-document numbers and refresh the
# "date issued" column for all 34 rows, then leave the F2:F34 range selected
# (mirrors what a user does after re-running the WG-number generator).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startNumber = 9271
$newDate = "10/19/2016"

for ($row = 1; $row -le 34; $row++) {
    $n = $startNumber + ($row - 1)
    $ws.Range("A$row").Value = "N$n"
    $ws.Range("F$row").Value = $newDate
}

# Matches the author's final selection state in the saved file.
$ws.Range("F2:F34").Select()
